$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (A3) held the wrong date (2013-10-07). Fix it to 2013-10-06.
$ws.Range("A3").Value = 41553

# Append a new row 4 for 2013-10-07 / 03:30:00, reusing row 3's formatting.
$ws.Range("A3:B3").Copy()
$ws.Range("A4:B4").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("A4").Value = 41554
$ws.Range("B4").Value = 0.14583333333333334

# Underline the "Quantidade de horas" column values (shared style covers B2:B4).
$ws.Range("B2:B4").Font.Underline = 2

# Match the selection left behind in the source workbook.
$ws.Range("B2:B4").Select()
